# Refresh computed market-price / profit columns (H:N) across the Leve
# profit sheets, per the scheduled Golem_Profits data-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 42: Eye of the Beholder
$ws.Range("H42").Value = 150.5
$ws.Range("I42").Value = 150.5
$ws.Range("K42").Value = 451.5
$ws.Range("M42").Value = -221.5

# Row 55: A Real Smooth Move
$ws.Range("H55").Value = 1614
$ws.Range("I55").Value = 2596.3333
$ws.Range("J55").Value = 631.6667
$ws.Range("K55").Value = 2596.3333
$ws.Range("L55").Value = 631.6667
$ws.Range("M55").Value = -2382.3333
$ws.Range("N55").Value = -1059.6667

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 2777.6
$ws.Range("I137").Value = 2699.6667
$ws.Range("K137").Value = 8099.000100000001
$ws.Range("M137").Value = -5549.000100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate
$ws.Range("H4").Value = 335.2
$ws.Range("I4").Value = 58.666668
$ws.Range("J4").Value = 750
$ws.Range("K4").Value = 58.666668
$ws.Range("L4").Value = 750
$ws.Range("M4").Value = 57.333332
$ws.Range("N4").Value = -982

# Row 58: Some Dragoons Have All the Luck
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1191.3334
$ws.Range("I122").Value = 1079.6
$ws.Range("K122").Value = 3238.8
$ws.Range("M122").Value = -788.7999999999997

# Row 123: The Armoire Is Open
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 3831.7
$ws.Range("I20").Value = 3349.5
$ws.Range("J20").Value = 4555
$ws.Range("K20").Value = 3349.5
$ws.Range("L20").Value = 4555
$ws.Range("M20").Value = -3102.5
$ws.Range("N20").Value = -5049

# Row 64: With Bearings Straight
$ws.Range("H64").Value = 450.7
$ws.Range("J64").Value = 626.3333
$ws.Range("L64").Value = 626.3333
$ws.Range("N64").Value = -1076.3333

# Row 67: Bearing the Brunt (L)
$ws.Range("H67").Value = 450.7
$ws.Range("J67").Value = 626.3333
$ws.Range("L67").Value = 626.3333
$ws.Range("N67").Value = -2186.3333

# Row 94: High Steal
$ws.Range("H94").Value = 2536.3635
$ws.Range("I94").Value = 1960
$ws.Range("K94").Value = 1960
$ws.Range("M94").Value = -1509

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 5428
$ws.Range("I99").Value = 6500
$ws.Range("K99").Value = 6500
$ws.Range("M99").Value = -5002

# Row 103: The Bigger the Blade
$ws.Range("H103").Value = 28073.889
$ws.Range("J103").Value = 28073.889
$ws.Range("L103").Value = 28073.889
$ws.Range("N103").Value = -30417.889

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 5432
$ws.Range("I134").Value = 5118.4
$ws.Range("K134").Value = 15355.2
$ws.Range("M134").Value = -12820.2

$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 3333499.8
$ws.Range("I4").Value = 499
$ws.Range("J4").Value = 5000000
$ws.Range("K4").Value = 499
$ws.Range("L4").Value = 5000000
$ws.Range("M4").Value = -387
$ws.Range("N4").Value = -5000224

# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1037
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 1271.8
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 1271.8
$ws.Range("M22").Value = -100
$ws.Range("N22").Value = -1971.8

# Row 31: Wall Not Found
$ws.Range("H31").Value = 5182.9443
$ws.Range("I31").Value = 4219.6
$ws.Range("J31").Value = 9999.666999999999
$ws.Range("K31").Value = 4219.6
$ws.Range("L31").Value = 9999.666999999999
$ws.Range("M31").Value = -3924.6
$ws.Range("N31").Value = -10589.667

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 5182.9443
$ws.Range("I34").Value = 4219.6
$ws.Range("J34").Value = 9999.666999999999
$ws.Range("K34").Value = 4219.6
$ws.Range("L34").Value = 9999.666999999999
$ws.Range("M34").Value = -4017.6
$ws.Range("N34").Value = -10403.667

# Row 35: Storm of Swords
$ws.Range("H35").Value = 4339.125
$ws.Range("J35").Value = 4300
$ws.Range("L35").Value = 4300
$ws.Range("N35").Value = -4888

# Row 44: Stay on Target
$ws.Range("H44").Value = 30071
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30071
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30071
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -30955

# Row 45: A Tree Grew in Gridania
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 3400
$ws.Range("I58").Value = 3400
$ws.Range("K58").Value = 3400
$ws.Range("M58").Value = -3197

# Row 92: Walk the Walk
$ws.Range("H92").Value = 68666.336
$ws.Range("J92").Value = 68666.336
$ws.Range("L92").Value = 68666.336
$ws.Range("N92").Value = -73658.336

# Row 93: Reeling for Rods
$ws.Range("H93").Value = 18347.75
$ws.Range("I93").Value = 18347.75
$ws.Range("K93").Value = 18347.75
$ws.Range("M93").Value = -16475.75

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 2466.6667
$ws.Range("I122").Value = 1700
$ws.Range("K122").Value = 5100
$ws.Range("M122").Value = -2650

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 2642
$ws.Range("I132").Value = 2012
$ws.Range("K132").Value = 6036
$ws.Range("M132").Value = -3506

# Row 136: Turali Quality
$ws.Range("H136").Value = 3400
$ws.Range("I136").Value = 3400
$ws.Range("K136").Value = 10200
$ws.Range("M136").Value = -7650

$ws = $wb.Worksheets.Item("CUL")
# Row 20: Omelette's Be Friends
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 26: A Grape Idea
$ws.Range("H26").Value = 97.666664
$ws.Range("I26").Value = 96.8
$ws.Range("K26").Value = 290.4
$ws.Range("M26").Value = -2.399999999999977

# Row 39: Bloody Good Tart, This
$ws.Range("H39").Value = 3141.5
$ws.Range("J39").Value = 3569.8
$ws.Range("L39").Value = 10709.4
$ws.Range("N39").Value = -11297.4

# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 230.11111
$ws.Range("J113").Value = 250
$ws.Range("L113").Value = 750
$ws.Range("N113").Value = -5090

# Row 122: Salt of the North
$ws.Range("H122").Value = 698.625
$ws.Range("I122").Value = 538.8
$ws.Range("J122").Value = 965
$ws.Range("K122").Value = 4849.2
$ws.Range("L122").Value = 8685
$ws.Range("M122").Value = -2399.2
$ws.Range("N122").Value = -13585

# Row 133: Friends Are Food
$ws.Range("H133").Value = 50014.5
$ws.Range("I133").Value = 50014.5
$ws.Range("K133").Value = 150043.5
$ws.Range("M133").Value = -144983.5

# Row 138: Bring Me Your Tacos
$ws.Range("H138").Value = 763.5714
$ws.Range("I138").Value = 763.5714
$ws.Range("K138").Value = 2290.7142
$ws.Range("M138").Value = 2849.2858

$ws = $wb.Worksheets.Item("GSM")
# Row 36: Keep the Change
$ws.Range("H36").Value = 14807.8
$ws.Range("I36").Value = 21179.666
$ws.Range("K36").Value = 21179.666
$ws.Range("M36").Value = -20694.666

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 1306.3334
$ws.Range("I97").Value = 1268.4
$ws.Range("J97").Value = 1496
$ws.Range("K97").Value = 1268.4
$ws.Range("L97").Value = 1496
$ws.Range("M97").Value = -772.4000000000001
$ws.Range("N97").Value = -2488

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 7866
$ws.Range("I126").Value = 9649
$ws.Range("J126").Value = 6974.5
$ws.Range("K126").Value = 28947
$ws.Range("L126").Value = 20923.5
$ws.Range("M126").Value = -26477
$ws.Range("N126").Value = -25863.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 2000
$ws.Range("J132").Value = 2000
$ws.Range("L132").Value = 6000
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore
$ws.Range("H16").Value = 791
$ws.Range("I16").Value = 882.7143
$ws.Range("J16").Value = 149
$ws.Range("K16").Value = 882.7143
$ws.Range("L16").Value = 149
$ws.Range("M16").Value = -712.7143
$ws.Range("N16").Value = -489

# Row 58: Handle with Care
$ws.Range("H58").Value = 32035
$ws.Range("I58").Value = 4105
$ws.Range("K58").Value = 4105
$ws.Range("M58").Value = -3845

# Row 111: Glove Me Tender
$ws.Range("H111").Value = 20000
$ws.Range("J111").Value = 20000
$ws.Range("L111").Value = 20000
$ws.Range("N111").Value = -28180

# Row 122: Hell on Leather
$ws.Range("H122").Value = 4897
$ws.Range("J122").Value = 4799
$ws.Range("L122").Value = 14397
$ws.Range("N122").Value = -19297

# Row 130: Generous Soles
$ws.Range("H130").Value = 41000
$ws.Range("J130").Value = 41000
$ws.Range("L130").Value = 41000
$ws.Range("N130").Value = -51040

$ws = $wb.Worksheets.Item("WVR")
# Row 13: Time for Acton
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 14: Hat in Hand
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Row 52: Party Animals
$ws.Range("H52").Value = 38990
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# Row 80: Healing with Flair
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()

# Row 83: Pants Fit for Battle (L)
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# Row 96: Skills on Display
$ws.Range("H96").Value = 24300
$ws.Range("I96").Value = 9000
$ws.Range("J96").Value = 28125
$ws.Range("K96").Value = 9000
$ws.Range("L96").Value = 28125
$ws.Range("M96").Value = -7627
$ws.Range("N96").Value = -30871

# Row 97: Getting a Leg Up
$ws.Range("H97").Value = 32499.834
$ws.Range("J97").Value = 32499.834
$ws.Range("L97").Value = 32499.834
$ws.Range("N97").Value = -34481.834

# Row 107: Flax Wax
$ws.Range("H107").Value = 1465.3334
$ws.Range("I107").Value = 1198
$ws.Range("K107").Value = 3594
$ws.Range("M107").Value = -1674
